$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Process Description paragraph: split the sentence and insert new text
#    "...specific language text file [through the text broker ]and returns..."
#    "...ing English [word][ if found. If not found it goes into the error file. ]"
# ---------------------------------------------------------------------------

# 1a) Insert "through the " + "text broker " right before
#     "and returns the translation of the correspond"
$rng = $d.Content
[void]$rng.Find.Execute("and returns the translation of the correspond", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)  # wdCollapseStart

$rng.InsertBefore("through the ")
$rng.Collapse(0)  # wdCollapseEnd
$rng.InsertBefore("text broker ")

# 1b) Replace the trailing space after "...ing English word" with
#     " if found. If not found it goes into the error file. "
$rng2 = $d.Content
[void]$rng2.Find.Execute("ing English word ", $true, $false, $false, $false, $false, $true, 1, $false, "ing English word if found. If not found it goes into the error file. ", 2)

# ---------------------------------------------------------------------------
# 2) Remove the block of empty "bold" paragraphs (and the "Pseudocode:"
#    paragraph) that used to sit between "APIs/Objects:" and the final blank
#    paragraph right before the section break.
# ---------------------------------------------------------------------------

$apiPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "APIs/Objects:") {
        $apiPara = $p
    }
}

if ($apiPara -ne $null) {
    $delStart = $apiPara.Range.End
    $scan = $apiPara.Next()
    $delEnd = $delStart
    $count = 0
    while ($scan -ne $null -and $count -lt 6) {
        $delEnd = $scan.Range.End
        $scan = $scan.Next()
        $count = $count + 1
    }
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
